# Updates cryptos list sheet with latest prices/volumes (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column holds numeric-looking price strings that must stay literal text
# (e.g. "69.070.14", "0.520", "1.00") - format the touched cells as Text first
# so Excel does not silently coerce them into Doubles and drop the formatting.
$dCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D9', 'D11', 'D12', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D36', 'D37', 'D39', 'D40', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50')
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '69.070.14'
$ws.Range('E2').Value = '  +1.29%  '
# Row 3
$ws.Range('D3').Value = '3.781.54'
$ws.Range('E3').Value = '  -0.53%  '
# Row 4
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.19%  '
# Row 5
$ws.Range('D5').Value = '628.75'
$ws.Range('E5').Value = '  +4.48%  '
# Row 6
$ws.Range('D6').Value = '164.15'
$ws.Range('E6').Value = '  -0.67%  '
# Row 7
$ws.Range('D7').Value = '3.777.32'
$ws.Range('E7').Value = '  -0.58%  '
# Row 8
$ws.Range('E8').Value = '  +0.00%  '
# Row 9
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').Value = '  +0.23%  '
# Row 10
$ws.Range('E10').Value = '  +0.80%  '
# Row 11
$ws.Range('D11').Value = '0.451'
$ws.Range('E11').Value = '  -0.14%  '
# Row 12
$ws.Range('D12').Value = '6.64'
$ws.Range('E12').Value = '  +2.41%  '
# Row 13
$ws.Range('E13').Value = '  -0.99%  '
# Row 14
$ws.Range('D14').Value = '35.41'
$ws.Range('E14').Value = '  -1.33%  '
# Row 15
$ws.Range('D15').Value = '4.411.06'
$ws.Range('E15').Value = '  -0.68%  '
# Row 16
$ws.Range('D16').Value = '3.711.52'
$ws.Range('E16').Value = '  -2.41%  '
# Row 17
$ws.Range('D17').Value = '68.947.36'
$ws.Range('E17').Value = '  +1.09%  '
# Row 18
$ws.Range('D18').Value = '17.93'
$ws.Range('E18').Value = '  -3.04%  '
# Row 19
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '7.09'
$ws.Range('E19').Value = '  -0.30%  '
# Row 20
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').Value = '0.113'
$ws.Range('E20').Value = '  -1.28%  '
# Row 21
$ws.Range('D21').Value = '467.88'
$ws.Range('E21').Value = '  +1.30%  '
# Row 22
$ws.Range('D22').Value = '9.65'
$ws.Range('E22').Value = '  -0.85%  '
# Row 23
$ws.Range('D23').Value = '0.703'
$ws.Range('E23').Value = '  +0.23%  '
# Row 24
$ws.Range('E24').Value = '  +0.29%  '
# Row 25
$ws.Range('D25').Value = '83.28'
$ws.Range('E25').Value = '  +0.15%  '
# Row 26
$ws.Range('D26').Value = '12.02'
$ws.Range('E26').Value = '  -0.48%  '
# Row 27
$ws.Range('D27').Value = '2.15'
$ws.Range('E27').Value = '  +1.47%  '
# Row 28
$ws.Range('E28').Value = '  +0.02%  '
# Row 29
$ws.Range('D29').Value = '10.00'
$ws.Range('E29').Value = '  -0.11%  '
# Row 30
$ws.Range('D30').Value = '3.924.08'
$ws.Range('E30').Value = '  -0.68%  '
# Row 31
$ws.Range('D31').Value = '2.67'
$ws.Range('E31').Value = '  +0.70%  '
# Row 32
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '7.27'
$ws.Range('E32').Value = '  -1.17%  '
# Row 33
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '2.21'
$ws.Range('E33').Value = '  -1.31%  '
# Row 34
$ws.Range('D34').Value = '28.96'
$ws.Range('E34').Value = '  -1.58%  '
# Row 35
$ws.Range('E35').Value = '  +0.06%  '
# Row 36
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = '8.99'
$ws.Range('E36').Value = '  -0.72%  '
# Row 37
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').Value = '3.719.28'
$ws.Range('E37').Value = '  -0.82%  '
# Row 38
$ws.Range('E38').Value = '  +2.47%  '
# Row 39
$ws.Range('D39').Value = '0.149'
$ws.Range('E39').Value = '  +7.06%  '
# Row 40
$ws.Range('D40').Value = '3.33'
$ws.Range('E40').Value = '  -0.01%  '
# Row 41
$ws.Range('D41').Value = '5.83'
$ws.Range('E41').Value = '  -0.39%  '
# Row 42
$ws.Range('D42').Value = '0.970'
$ws.Range('E42').Value = '  -1.95%  '
# Row 43
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.18%  '
# Row 44
$ws.Range('E44').Value = '  +0.03%  '
# Row 45
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '0.298'
$ws.Range('E45').Value = '  -0.93%  '
# Row 46
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = '153.40'
$ws.Range('E46').Value = '  +1.24%  '
# Row 47
$ws.Range('D47').Value = '46.87'
$ws.Range('E47').Value = '  -1.43%  '
# Row 48
$ws.Range('D48').Value = '1.92'
$ws.Range('E48').Value = '  +2.38%  '
# Row 49
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').Value = '42.66'
$ws.Range('E49').Value = '  -1.31%  '
# Row 50
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '8.41'
$ws.Range('E50').Value = '  +0.33%  '
# Row 51
$ws.Range('E51').Value = '  +1.66%  '
